$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2429.3
$ws.Range("I9").Value = 2889.5715
$ws.Range("J9").Value = 1355.3334
$ws.Range("K9").Value = 2889.5715
$ws.Range("L9").Value = 1355.3334
$ws.Range("M9").Value = -2720.5715
$ws.Range("N9").Value = -1693.3334

$ws.Range("H32").Value = 15651.037
$ws.Range("I32").Value = 17912
$ws.Range("J32").Value = 13551.571
$ws.Range("K32").Value = 17912
$ws.Range("L32").Value = 13551.571
$ws.Range("M32").Value = -17586
$ws.Range("N32").Value = -14203.571

$ws.Range("H88").Value = 2313.1428
$ws.Range("I88").Value = 799
$ws.Range("J88").Value = 2429.6155
$ws.Range("K88").Value = 799
$ws.Range("L88").Value = 2429.6155
$ws.Range("M88").Value = -393
$ws.Range("N88").Value = -3241.6155

$ws.Range("H91").Value = 2313.1428
$ws.Range("I91").Value = 799
$ws.Range("J91").Value = 2429.6155
$ws.Range("K91").Value = 799
$ws.Range("L91").Value = 2429.6155
$ws.Range("M91").Value = 605
$ws.Range("N91").Value = -5237.6155

$ws.Range("H116").Value = 4837847
$ws.Range("I116").Value = 5853972.5
$ws.Range("J116").Value = 11249.25
$ws.Range("K116").Value = 5853972.5
$ws.Range("L116").Value = 11249.25
$ws.Range("M116").Value = -5850530.5
$ws.Range("N116").Value = -18133.25

$ws.Range("H137").Value = 6285.4814
$ws.Range("I137").Value = 3668.8667
$ws.Range("J137").Value = 9556.25
$ws.Range("K137").Value = 11006.6001
$ws.Range("L137").Value = 28668.75
$ws.Range("M137").Value = -8456.6001
$ws.Range("N137").Value = -33768.75

$ws.Range("H138").Value = 5596.255
$ws.Range("I138").Value = 2633.6667
$ws.Range("J138").Value = 6507.8203
$ws.Range("K138").Value = 7901.000100000001
$ws.Range("L138").Value = 19523.4609
$ws.Range("M138").Value = -2761.000100000001
$ws.Range("N138").Value = -29803.4609

$ws.Range("H141").Value = 1665.0667
$ws.Range("I141").Value = 1398.8
$ws.Range("J141").Value = 2197.6
$ws.Range("K141").Value = 4196.4
$ws.Range("L141").Value = 6592.799999999999
$ws.Range("M141").Value = 983.6000000000004
$ws.Range("N141").Value = -16952.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6449.6665
$ws.Range("I2").Value = 3819.4119
$ws.Range("J2").Value = 17628.25
$ws.Range("K2").Value = 3819.4119
$ws.Range("L2").Value = 17628.25
$ws.Range("M2").Value = -3706.4119

$ws.Range("H19").Value = 16803334
$ws.Range("I19").Value = 50000000
$ws.Range("J19").Value = 205000
$ws.Range("K19").Value = 50000000
$ws.Range("L19").Value = 205000
$ws.Range("M19").Value = -49999771
$ws.Range("N19").Value = -205458

$ws.Range("H32").Value = 463.73
$ws.Range("I32").Value = 453.06122
$ws.Range("J32").Value = 986.5
$ws.Range("K32").Value = 453.06122
$ws.Range("L32").Value = 986.5
$ws.Range("M32").Value = -166.06122
$ws.Range("N32").Value = -1560.5

$ws.Range("H88").Value = 7524.727
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 7524.727
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 7524.727
$ws.Range("N88").Value = -8336.726999999999

$ws.Range("H91").Value = 7524.727
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 7524.727
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 7524.727
$ws.Range("N91").Value = -10332.727

$ws.Range("H116").Value = 6449.6665
$ws.Range("I116").Value = 3819.4119
$ws.Range("J116").Value = 17628.25
$ws.Range("K116").Value = 3819.4119
$ws.Range("L116").Value = 17628.25
$ws.Range("M116").Value = -1525.4119

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6449.6665
$ws.Range("I3").Value = 3819.4119
$ws.Range("J3").Value = 17628.25
$ws.Range("K3").Value = 3819.4119
$ws.Range("L3").Value = 17628.25
$ws.Range("M3").Value = -3705.4119

$ws.Range("H86").Value = 5818.8335
$ws.Range("I86").Value = 4735.4
$ws.Range("J86").Value = 8281.182000000001
$ws.Range("K86").Value = 4735.4
$ws.Range("L86").Value = 8281.182000000001
$ws.Range("M86").Value = -3612.4

$ws.Range("H89").Value = 5818.8335
$ws.Range("I89").Value = 4735.4
$ws.Range("J89").Value = 8281.182000000001
$ws.Range("K89").Value = 23677
$ws.Range("L89").Value = 41405.91
$ws.Range("M89").Value = -18061

$ws.Range("H134").Value = 1229721.6
$ws.Range("I134").Value = 1330232.4
$ws.Range("J134").Value = 23592.666
$ws.Range("K134").Value = 3990697.2
$ws.Range("L134").Value = 70777.99800000001
$ws.Range("M134").Value = -3988162.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9500
$ws.Range("I4").Value = 9500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 9500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -9388

$ws.Range("H107").Value = 677.4211
$ws.Range("I107").Value = 590.53845
$ws.Range("J107").Value = 865.6667
$ws.Range("K107").Value = 590.53845
$ws.Range("L107").Value = 865.6667
$ws.Range("M107").Value = 1329.46155
$ws.Range("N107").Value = -4705.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 131372.56
$ws.Range("I122").Value = 310.375
$ws.Range("J122").Value = 165195.06
$ws.Range("K122").Value = 2793.375
$ws.Range("L122").Value = 1486755.54
$ws.Range("M122").Value = -343.375
$ws.Range("N122").Value = -1491655.54

$ws.Range("H127").Value = 3866
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 3866
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 11598
$ws.Range("N127").Value = -21518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8692.362999999999
$ws.Range("I70").Value = 5885.643
$ws.Range("J70").Value = 13604.125
$ws.Range("K70").Value = 5885.643
$ws.Range("L70").Value = 13604.125
$ws.Range("M70").Value = -5615.643
$ws.Range("N70").Value = -14144.125

$ws.Range("H73").Value = 8692.362999999999
$ws.Range("I73").Value = 5885.643
$ws.Range("J73").Value = 13604.125
$ws.Range("K73").Value = 5885.643
$ws.Range("L73").Value = 13604.125
$ws.Range("M73").Value = -4949.643
$ws.Range("N73").Value = -15476.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1010
$ws.Range("I30").Value = 1010
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1010
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -902

$ws.Range("H46").Value = 125003120
$ws.Range("I46").Value = 9999
$ws.Range("J46").Value = 166667500
$ws.Range("K46").Value = 9999
$ws.Range("L46").Value = 166667500
$ws.Range("M46").Value = -9811
$ws.Range("N46").Value = -166667876

$ws.Range("H132").Value = 6523.8335
$ws.Range("I132").Value = 6464.316
$ws.Range("J132").Value = 6750
$ws.Range("K132").Value = 19392.948
$ws.Range("L132").Value = 20250
$ws.Range("M132").Value = -16862.948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 37499.5
$ws.Range("I49").Value = 49999
$ws.Range("J49").Value = 25000
$ws.Range("K49").Value = 49999
$ws.Range("L49").Value = 25000
$ws.Range("M49").Value = -49769
$ws.Range("N49").Value = -25460

$ws.Range("H107").Value = 4517.8
$ws.Range("I107").Value = 4517.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 13553.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -11633.4
$ws.Range("N107").ClearContents()

$ws.Range("H110").Value = 74000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 74000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 74000
$ws.Range("N110").Value = -82180

$ws.Range("H113").Value = 4904250.5
$ws.Range("I113").Value = 5954352.5
$ws.Range("J113").Value = 3773.3333
$ws.Range("K113").Value = 17863057.5
$ws.Range("L113").Value = 11319.9999
$ws.Range("M113").Value = -17860887.5
$ws.Range("N113").Value = -15659.9999

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H122").Value = 5833
$ws.Range("I122").Value = 5833
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 17499
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15049
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 19750
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 19750
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 59250
$ws.Range("N126").Value = -64190
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 9665.77
$ws.Range("I132").Value = 9543
$ws.Range("J132").Value = 9942
$ws.Range("K132").Value = 28629
$ws.Range("L132").Value = 29826
$ws.Range("M132").Value = -26099
$ws.Range("N132").Value = -34886

$ws.Range("H136").Value = 9264728
$ws.Range("I136").Value = 10420747
$ws.Range("J136").Value = 16582.834
$ws.Range("K136").Value = 31262241
$ws.Range("L136").Value = 49748.50199999999
$ws.Range("M136").Value = -31259691

$ws.Range("H137").Value = 80000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 80000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200
